$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.86
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 5.1
$ws.Range("J2").Value = 3.4
$ws.Range("N2").Value = 3.4
$ws.Range("Q2").Value = 1.91
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 1.96
$ws.Range("V2").Value = 1.25
$ws.Range("Y2").Value = 18.5
$ws.Range("AF2").Value = 14.5
$ws.Range("AO2").Value = 90
$ws.Range("L3").Value = 1.39
$ws.Range("O3").Value = 1.31
$ws.Range("Y3").Value = 17
$ws.Range("G4").Value = 2.78
$ws.Range("J4").Value = 2.84
$ws.Range("L4").Value = 1.48
$ws.Range("Q4").Value = 2.08
$ws.Range("V4").Value = 1.34
$ws.Range("W4").Value = 1.59
$ws.Range("F5").Value = 1.17
$ws.Range("G5").Value = 1.23
$ws.Range("I5").Value = 40
$ws.Range("K5").Value = 8.8
$ws.Range("AF5").Value = 7.4
$ws.Range("AH5").Value = 75
$ws.Range("AJ5").Value = 1000
$ws.Range("AN5").Value = 4.6
$ws.Range("G6").Value = 1.84
$ws.Range("T7").Value = 2.22
$ws.Range("U7").Value = 1.8
$ws.Range("G8").Value = 2.28
$ws.Range("H8").Value = 3.2
$ws.Range("L8").Value = 1.25
$ws.Range("AA8").Value = 70
$ws.Range("AB8").Value = 16.5
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 18.5
$ws.Range("AI8").Value = 46
$ws.Range("AK8").Value = 25
$ws.Range("AL8").Value = 36
$ws.Range("AM8").Value = 75
$ws.Range("AO8").Value = 32
$ws.Range("F9").Value = 2.26
$ws.Range("K9").Value = 4.6
$ws.Range("V9").Value = 1.5
$ws.Range("L10").Value = 1.19
$ws.Range("N10").Value = 7
$ws.Range("Q10").Value = 1.4
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 1.48
$ws.Range("W10").Value = 1.79
$ws.Range("F11").Value = 1.87
$ws.Range("G11").Value = 1.99
$ws.Range("H11").Value = 3.8
$ws.Range("J11").Value = 4.1
$ws.Range("K11").Value = 4.6
$ws.Range("P11").Value = 2.88
$ws.Range("Q11").Value = 1.43
$ws.Range("S11").Value = 2.08
$ws.Range("U11").Value = 2.74
$ws.Range("W11").Value = 2
$ws.Range("X11").Value = 40
$ws.Range("Z11").Value = 44
$ws.Range("AC11").Value = 14
$ws.Range("AD11").Value = 21
$ws.Range("AH11").Value = 18
$ws.Range("L12").Value = 1.43
$ws.Range("S12").Value = 3.75
$ws.Range("AC12").Value = 9.2
$ws.Range("AK12").Value = 980
$ws.Range("I13").Value = 2.84
$ws.Range("R13").Value = 1.23
$ws.Range("T13").Value = 1.95
$ws.Range("V13").Value = 1.54
$ws.Range("AC13").Value = 7.2
$ws.Range("G14").Value = 1.27
$ws.Range("H14").Value = 15
$ws.Range("Q14").Value = 1.69
$ws.Range("T14").Value = 2.4
$ws.Range("U14").Value = 1.67
$ws.Range("W14").Value = 4.7
$ws.Range("AB14").Value = 8.4
$ws.Range("AJ14").Value = 8.8
$ws.Range("H15").Value = 2.8
$ws.Range("I15").Value = 2.84
$ws.Range("J15").Value = 3.35
$ws.Range("V15").Value = 1.54
$ws.Range("X15").Value = 12.5
